# "updating prices at 10:57:59" -- append a new price snapshot as row 22.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 22
$prevRow = $newRow - 1

$values = @(32.2, 8.4, 26.3, 6.4, 24.3, 12.6, 12.1, 10.8, 7.4, 15.9)

for ($col = 1; $col -le $values.Length; $col++) {
    $cell = $ws.Cells.Item($newRow, $col)
    $cell.Value = $values[$col - 1]
}

# Timestamp column (K) - carry over the same date/time number format used
# by the rows above it.
$tsCell = $ws.Cells.Item($newRow, 11)
$tsCell.Value = 44042.91666666666
$tsCell.NumberFormat = $ws.Cells.Item($prevRow, 11).NumberFormat
